# Regenerate the "K" column (column G) values for the save-data sheet.
# Per the commit message, column G ("K") values were recomputed (std/mean
# recalculated, s_vals written), so we update each data row's G cell with
# its newly computed value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => new K value, as produced by the regenerated calculation.
$kValues = @{
    2  = 0
    3  = 2
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 1
    9  = 1
    10 = 1
    11 = 3
    12 = 2
    13 = 1
    14 = 1
    15 = 1
    16 = 1
    17 = 2
    18 = 2
    19 = 2
    20 = 1
    21 = 0
    22 = 1
    23 = 1
    24 = 2
    25 = 1
    26 = 2
    27 = 3
    28 = 0
    29 = 3
    30 = 0
    31 = 1
    32 = 2
    33 = 3
    34 = 0
    35 = 1
    36 = 1
    37 = 1
    38 = 2
    39 = 3
    40 = 2
    41 = 0
    42 = 0
    43 = 1
    44 = 1
    45 = 0
    46 = 3
    47 = 2
    48 = 0
    49 = 1
    50 = 1
    51 = 1
    52 = 2
    53 = 1
    54 = 2
    55 = 3
    56 = 1
    58 = 1
    59 = 1
    60 = 0
    61 = 1
    62 = 0
    63 = 2
    64 = 0
    65 = 3
    66 = 1
    67 = 0
    68 = 2
    69 = 2
    70 = 0
    71 = 3
    72 = 1
    73 = 1
    74 = 2
    75 = 1
    76 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
